# Insert a new row at position 14, pushing the existing rows 14-118 down
# to 15-119 (matches the dimension change from A1:R118 to A1:R119), then
# populate the newly-inserted row with the new "Arveja Verde" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("14").Insert()

$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 44970
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 100112022
$ws.Range("G14").Value = "Arveja Verde"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 70
$ws.Range("K14").Value = 25000
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = 27143
$ws.Range("N14").Value = "$/saco 25 kilos"
$ws.Range("O14").Value = "Región de La Araucanía"
$ws.Range("P14").Value = 1086
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
